$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the latest crypto snapshot: updated prices / 1h volume %
# deltas, plus the Filecoin/PancakeSwap ranking swap (rows 32-33).
# Price column (D) holds plain-text numeric-looking strings (e.g.
# "11.00", "225.88") in the source data; force text format before
# assigning so Excel does not silently reinterpret them as numbers
# and strip significant trailing zeros / separators.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.364.09"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.13"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.88"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.62"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0945"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.047.44"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.793.15"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.00"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.369.93"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.27"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "244.17"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.21"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.15"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.81"
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  +6.00%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.79"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.58"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.400.83"
$ws.Range("E36").Value = "  -3.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.674"
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.47"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.934"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.84"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0527"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.99"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.947.84"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.61"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  -3.05%  "
